$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 93.24418900000001
$ws.Range("H2").Value = 279.732567
$ws.Range("I2").Value = 0.5033448418000582
$ws.Range("J2").Value = 0.5033448418000582
$ws.Range("M2").Value = 184.1023456666667
$ws.Range("N2").Value = 552.307037
$ws.Range("O2").Value = 0.9813423747591566
$ws.Range("P2").Value = 0.9813423747591565
$ws.Range("Q2").Value = 17166.473914686
$ws.Range("R2").Value = 154498.265232174
$ws.Range("S2").Value = 0.4939536223748411
$ws.Range("T2").Value = 0.4939536223748411
$ws.Range("G3").Value = 93.24418900000001
$ws.Range("H3").Value = 279.732567
$ws.Range("I3").Value = 0.5033448418000582
$ws.Range("J3").Value = 0.5033448418000582
$ws.Range("M3").Value = 0.4321196666666667
$ws.Range("N3").Value = 1.296359
$ws.Range("O3").Value = 0.002303378255889225
$ws.Range("P3").Value = 0.002303378255889224
$ws.Range("Q3").Value = 40.29264786928367
$ws.Range("R3").Value = 362.633830823553
$ws.Range("S3").Value = 0.001159393563816256
$ws.Range("T3").Value = 0.001159393563816256
$ws.Range("G4").Value = 93.24418900000001
$ws.Range("H4").Value = 279.732567
$ws.Range("I4").Value = 0.5033448418000582
$ws.Range("J4").Value = 0.5033448418000582
$ws.Range("M4").Value = 1.367901
$ws.Range("N4").Value = 4.103703
$ws.Range("O4").Value = 0.007291483500193526
$ws.Range("P4").Value = 0.007291483500193526
$ws.Range("Q4").Value = 127.548819377289
$ws.Range("R4").Value = 1147.939374395601
$ws.Range("S4").Value = 0.003670130608892645
$ws.Range("T4").Value = 0.003670130608892645
$ws.Range("G5").Value = 93.24418900000001
$ws.Range("H5").Value = 279.732567
$ws.Range("I5").Value = 0.5033448418000582
$ws.Range("J5").Value = 0.5033448418000582
$ws.Range("M5").Value = 1.700197666666667
$ws.Range("N5").Value = 5.100593
$ws.Range("O5").Value = 0.009062763484760617
$ws.Range("P5").Value = 0.009062763484760615
$ws.Range("Q5").Value = 158.5335525680257
$ws.Range("R5").Value = 1426.801973112231
$ws.Range("S5").Value = 0.004561695252508177
$ws.Range("T5").Value = 0.004561695252508176
$ws.Range("G6").Value = 31.60427266666667
$ws.Range("H6").Value = 94.81281800000001
$ws.Range("I6").Value = 0.1706041716509459
$ws.Range("J6").Value = 0.1706041716509459
$ws.Range("M6").Value = 184.1023456666667
$ws.Range("N6").Value = 552.307037
$ws.Range("O6").Value = 0.9813423747591566
$ws.Range("P6").Value = 0.9813423747591565
$ws.Range("Q6").Value = 5818.420731022253
$ws.Range("R6").Value = 52365.78657920027
$ws.Range("S6").Value = 0.167421102951758
$ws.Range("T6").Value = 0.167421102951758
$ws.Range("G7").Value = 31.60427266666667
$ws.Range("H7").Value = 94.81281800000001
$ws.Range("I7").Value = 0.1706041716509459
$ws.Range("J7").Value = 0.1706041716509459
$ws.Range("M7").Value = 0.4321196666666667
$ws.Range("N7").Value = 1.296359
$ws.Range("O7").Value = 0.002303378255889225
$ws.Range("P7").Value = 0.002303378255889224
$ws.Range("Q7").Value = 13.65682776996245
$ws.Range("R7").Value = 122.911449929662
$ws.Range("S7").Value = 0.0003929659393447816
$ws.Range("T7").Value = 0.0003929659393447816
$ws.Range("G8").Value = 31.60427266666667
$ws.Range("H8").Value = 94.81281800000001
$ws.Range("I8").Value = 0.1706041716509459
$ws.Range("J8").Value = 0.1706041716509459
$ws.Range("M8").Value = 1.367901
$ws.Range("N8").Value = 4.103703
$ws.Range("O8").Value = 0.007291483500193526
$ws.Range("P8").Value = 0.007291483500193526
$ws.Range("Q8").Value = 43.23151618500601
$ws.Range("R8").Value = 389.0836456650541
$ws.Range("S8").Value = 0.001243957502657056
$ws.Range("T8").Value = 0.001243957502657056
$ws.Range("G9").Value = 31.60427266666667
$ws.Range("H9").Value = 94.81281800000001
$ws.Range("I9").Value = 0.1706041716509459
$ws.Range("J9").Value = 0.1706041716509459
$ws.Range("M9").Value = 1.700197666666667
$ws.Range("N9").Value = 5.100593
$ws.Range("O9").Value = 0.009062763484760617
$ws.Range("P9").Value = 0.009062763484760615
$ws.Range("Q9").Value = 53.73351064456379
$ws.Range("R9").Value = 483.601595801074
$ws.Range("S9").Value = 0.001546145257186025
$ws.Range("T9").Value = 0.001546145257186025
$ws.Range("G10").Value = 60.37827433333333
$ws.Range("H10").Value = 181.134823
$ws.Range("I10").Value = 0.3259301546659619
$ws.Range("J10").Value = 0.3259301546659619
$ws.Range("M10").Value = 184.1023456666667
$ws.Range("N10").Value = 552.307037
$ws.Range("O10").Value = 0.9813423747591566
$ws.Range("P10").Value = 0.9813423747591565
$ws.Range("Q10").Value = 11115.78193207216
$ws.Range("R10").Value = 100042.0373886494
$ws.Range("S10").Value = 0.3198490719855143
$ws.Range("T10").Value = 0.3198490719855142
$ws.Range("G11").Value = 60.37827433333333
$ws.Range("H11").Value = 181.134823
$ws.Range("I11").Value = 0.3259301546659619
$ws.Range("J11").Value = 0.3259301546659619
$ws.Range("M11").Value = 0.4321196666666667
$ws.Range("N11").Value = 1.296359
$ws.Range("O11").Value = 0.002303378255889225
$ws.Range("P11").Value = 0.002303378255889224
$ws.Range("Q11").Value = 26.09063977882856
$ws.Range("R11").Value = 234.815758009457
$ws.Range("S11").Value = 0.0007507404311961885
$ws.Range("T11").Value = 0.0007507404311961884
$ws.Range("G12").Value = 60.37827433333333
$ws.Range("H12").Value = 181.134823
$ws.Range("I12").Value = 0.3259301546659619
$ws.Range("J12").Value = 0.3259301546659619
$ws.Range("M12").Value = 1.367901
$ws.Range("N12").Value = 4.103703
$ws.Range("O12").Value = 0.007291483500193526
$ws.Range("P12").Value = 0.007291483500193526
$ws.Range("Q12").Value = 82.591501838841
$ws.Range("R12").Value = 743.323516549569
$ws.Range("S12").Value = 0.002376514344962385
$ws.Range("T12").Value = 0.002376514344962385
$ws.Range("G13").Value = 60.37827433333333
$ws.Range("H13").Value = 181.134823
$ws.Range("I13").Value = 0.3259301546659619
$ws.Range("J13").Value = 0.3259301546659619
$ws.Range("M13").Value = 1.700197666666667
$ws.Range("N13").Value = 5.100593
$ws.Range("O13").Value = 0.009062763484760617
$ws.Range("P13").Value = 0.009062763484760615
$ws.Range("Q13").Value = 102.6550011388932
$ws.Range("R13").Value = 923.8950102500389
$ws.Range("S13").Value = 0.002953827904289059
$ws.Range("T13").Value = 0.002953827904289059
$ws.Range("G14").Value = 0.022384
$ws.Range("H14").Value = 0.067152
$ws.Range("I14").Value = 0.0001208318830340463
$ws.Range("J14").Value = 0.0001208318830340463
$ws.Range("M14").Value = 184.1023456666667
$ws.Range("N14").Value = 552.307037
$ws.Range("O14").Value = 0.9813423747591566
$ws.Range("P14").Value = 0.9813423747591565
$ws.Range("Q14").Value = 4.120946905402668
$ws.Range("R14").Value = 37.088522148624
$ws.Range("S14").Value = 0.0001185774470432516
$ws.Range("T14").Value = 0.0001185774470432516
$ws.Range("G15").Value = 0.022384
$ws.Range("H15").Value = 0.067152
$ws.Range("I15").Value = 0.0001208318830340463
$ws.Range("J15").Value = 0.0001208318830340463
$ws.Range("M15").Value = 0.4321196666666667
$ws.Range("N15").Value = 1.296359
$ws.Range("O15").Value = 0.002303378255889225
$ws.Range("P15").Value = 0.002303378255889224
$ws.Range("Q15").Value = 0.009672566618666667
$ws.Range("R15").Value = 0.087053099568
$ws.Range("S15").Value = [double]"2.783215319987723E-07"
$ws.Range("T15").Value = [double]"2.783215319987722E-07"
$ws.Range("G16").Value = 0.022384
$ws.Range("H16").Value = 0.067152
$ws.Range("I16").Value = 0.0001208318830340463
$ws.Range("J16").Value = 0.0001208318830340463
$ws.Range("M16").Value = 1.367901
$ws.Range("N16").Value = 4.103703
$ws.Range("O16").Value = 0.007291483500193526
$ws.Range("P16").Value = 0.007291483500193526
$ws.Range("Q16").Value = 0.030619095984
$ws.Range("R16").Value = 0.275571863856
$ws.Range("S16").Value = [double]"8.810436814400624E-07"
$ws.Range("T16").Value = [double]"8.810436814400624E-07"
$ws.Range("G17").Value = 0.022384
$ws.Range("H17").Value = 0.067152
$ws.Range("I17").Value = 0.0001208318830340463
$ws.Range("J17").Value = 0.0001208318830340463
$ws.Range("M17").Value = 1.700197666666667
$ws.Range("N17").Value = 5.100593
$ws.Range("O17").Value = 0.009062763484760617
$ws.Range("P17").Value = 0.009062763484760615
$ws.Range("Q17").Value = 0.03805722457066667
$ws.Range("R17").Value = 0.342515021136
$ws.Range("S17").Value = [double]"1.09507077735582E-06"
$ws.Range("T17").Value = [double]"1.09507077735582E-06"

Write-Output "Applied 192 cell updates"
